# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the run now produced a
# new GUID-named markdown/xliff file pair and new handoff timestamps.
# Update the three report sheets (Overview, zh-cn, de-de) accordingly,
# including refreshing the display text of the existing "e2e\<id>.md"
# hyperlinks so they keep pointing at the same target URL but show the
# new file name.

$wb = $excel.ActiveWorkbook

$oldId = "18db6903-19e1-4f36-876f-f9896f2a5076"
$newId = "42ac6cd6-cf51-410a-88f8-4a761a316219"

$oldHash = "1b46799161dc36333e8369ade67fa345816fff60"
$newHash = "d1a9fbbff224a28d2e68216e4c57beeaedb644bd"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e6169f5fa2cc675e903f082959035eeae5d41ae5/e2e/$oldId.md"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("G2").Value = "2016-08-20 12:59:41"
# B2 carries the hyperlink; refresh its display text while keeping the
# same external target address.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newId.md")

# --- zh-cn sheet ----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 12:59:36"
# A2 carries the hyperlink; refresh its display text the same way.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newId.md")

# --- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-20 12:59:41"
# A2 carries the hyperlink; refresh its display text the same way.
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newId.md")
